# Add "ECA" and "OHI" 2021 "all" rows to the "regions" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("regions")

# Row 20: ECA, 2021, all
$ws.Range("A20").Value = "ECA"
$ws.Range("B20").Value = 2021
$ws.Range("B20").NumberFormat = "0"
$ws.Range("C20").Value = "all"

# Row 21: OHI, 2021, all
$ws.Range("A21").Value = "OHI"
$ws.Range("B21").Value = 2021
$ws.Range("B21").NumberFormat = "0"
$ws.Range("C21").Value = "all"

# Match the selection left by the author after adding the new rows
$ws.Range("A23").Select()
